$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col2"
$ws.Range("C1").Value = "col3"
$ws.Range("D1").Value = "col4"
$ws.Range("E1").Value = "col5"
$ws.Range("F1").Value = "col6"

$range = $ws.Range("A1:F1")
$range.Font.Bold = $true
$range.HorizontalAlignment = -4108
$range.Select()

